$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A30").Value = "Trial 29"
$ws.Range("B30").Value = 88.65517687797546
$ws.Range("C30").Value = 0.721221923828125
$ws.Range("D30").Value = 0.1725394725799561
$ws.Range("E30").Value = 0.2868828773498535
$ws.Range("F30").Value = 0.1107616424560547
$ws.Range("G30").Value = 0.1794347763061523
$ws.Range("H30").Value = 0.01595711708068848
